# Revert "Added instructions for working on a github project"
#
# Removes the "Keeping your project up-to-date" section (a heading
# paragraph plus the five numbered/bulleted steps that follow it, one of
# which carries the hidden "_GoBack" bookmark) that had been appended
# after the "Step-by-step guide..." bullet, restoring the document to
# its pre-edit state: a single trailing empty paragraph that still hosts
# the "_GoBack" bookmark.

$d = $word.ActiveDocument

# Locate the paragraph immediately after "Step-by-step guide..." (the
# first paragraph of the block being removed) and the paragraph that
# carries the "Press 'Sync' ..." text (the last paragraph of the block,
# which also holds the bookmark we must preserve).
$removeFromIndex = -1
$removeThroughIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Step-by-step guide at desktop.github.com*") {
        $removeFromIndex = $i + 1
    }
    if ($t -like "*Press*Sync*push those changes*") {
        $removeThroughIndex = $i
    }
}

if ($removeFromIndex -ne -1 -and $removeThroughIndex -ne -1 -and $removeThroughIndex -ge $removeFromIndex) {

    # Delete everything from the start of the first unwanted paragraph
    # through the start of the paragraph that follows the last unwanted
    # paragraph (i.e. the whole block, including every paragraph mark
    # inside it, but stopping right before the final trailing empty
    # paragraph so that paragraph - and the document's last paragraph
    # mark - survive untouched).
    $blockStart = $d.Paragraphs.Item($removeFromIndex).Range.Start
    $blockEnd = $d.Paragraphs.Item($removeThroughIndex + 1).Range.Start
    $d.Range($blockStart, $blockEnd).Delete()

    # Re-create the "_GoBack" bookmark (previously inside the deleted
    # "Press 'Sync' ..." paragraph) on the now-empty final paragraph, so
    # the document ends exactly as it did before the edit.
    $lastIndex = $d.Paragraphs.Count
    $lastRange = $d.Paragraphs.Item($lastIndex).Range
    $d.Bookmarks.Add("_GoBack", $lastRange)
}
